$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark row 7 (Nacho) as sent
$ws.Range("C7").Value = "y"

# Add row 9: BuffyGirl
$ws.Range("A9").Value = "BuffyGirl"
$ws.Range("B9").Value = "JC Kovacs`n12349 Metric Blvd #1330`nAustin, TX`n78758"
$ws.Range("B9").Style = $ws.Range("B5").Style
$ws.Range("C9").Value = "y"

# Add row 10: King Darkness
$ws.Range("A10").Value = "King Darkness"
$ws.Range("B10").Value = "Chase Valdez`n2418 East Highway 66`nPMB 539`nGallup NM 87301"
$ws.Range("B10").Style = $ws.Range("B5").Style
$ws.Range("C10").Value = "y"

# Update the view: scroll so A5 is the top-left visible cell, and select C7
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C7").Select()
